$wb = $excel.ActiveWorkbook

# Update "想去人数" (interested count) figures on the sheets that contain
# the exhibition data: "展览" and "全部类型" (they hold identical rows).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 1516
    $ws.Range("F9").Value = 330
}
